# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to per-sheet market-board profit tables
# (currentAveragePrice / LevePrice / LeveProfit columns), per commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 49142.934
$ws.Range("J17").Value = 49142.934
$ws.Range("L17").Value = 147428.802
$ws.Range("N17").Value = -147764.802
$ws.Range("H74").Value = 7598.353
$ws.Range("I74").Value = 7869.7144
$ws.Range("J74").Value = 6332
$ws.Range("K74").Value = 7869.7144
$ws.Range("L74").Value = 6332
$ws.Range("M74").Value = -6933.7144
$ws.Range("N74").Value = -8204
$ws.Range("H77").Value = 7598.353
$ws.Range("I77").Value = 7869.7144
$ws.Range("J77").Value = 6332
$ws.Range("K77").Value = 39348.572
$ws.Range("L77").Value = 31660
$ws.Range("M77").Value = -34668.572
$ws.Range("N77").Value = -41020
$ws.Range("H94").Value = 832.46155
$ws.Range("I94").Value = 832.46155
$ws.Range("K94").Value = 832.46155
$ws.Range("M94").Value = -381.46155
$ws.Range("H127").Value = 2529.4
$ws.Range("I127").Value = 2529.4
$ws.Range("K127").Value = 7588.200000000001
$ws.Range("M127").Value = -2628.200000000001
$ws.Range("H135").Value = 1094.3334
$ws.Range("I135").Value = 1094.3334
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 9849.000599999999
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 1731.6428
$ws.Range("I137").Value = 1457.6666
$ws.Range("K137").Value = 4372.9998
$ws.Range("M137").Value = -1822.9998
$ws.Range("H141").Value = 4039.6
$ws.Range("I141").Value = 4349.25
$ws.Range("J141").Value = 3833.1667
$ws.Range("K141").Value = 13047.75
$ws.Range("L141").Value = 11499.5001
$ws.Range("M141").Value = -7867.75
$ws.Range("N141").Value = -21859.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 169098.5
$ws.Range("I45").Value = 252597.75
$ws.Range("K45").Value = 252597.75
$ws.Range("M45").Value = -252220.75
$ws.Range("H74").Value = 2898.0967
$ws.Range("I74").Value = 1239.5122
$ws.Range("K74").Value = 1239.5122
$ws.Range("M74").Value = -365.5121999999999
$ws.Range("H77").Value = 2898.0967
$ws.Range("I77").Value = 1239.5122
$ws.Range("K77").Value = 6197.561
$ws.Range("M77").Value = -1829.561
$ws.Range("H92").Value = 59991.668
$ws.Range("J92").Value = 59991.668
$ws.Range("L92").Value = 59991.668
$ws.Range("N92").Value = -64983.668
$ws.Range("H110").Value = 562.5
$ws.Range("I110").Value = 270
$ws.Range("J110").Value = 855
$ws.Range("K110").Value = 270
$ws.Range("L110").Value = 855
$ws.Range("M110").Value = 1775
$ws.Range("N110").Value = -4945
$ws.Range("H122").Value = 4298.8335
$ws.Range("I122").Value = 4401
$ws.Range("K122").Value = 13203
$ws.Range("M122").Value = -10753
$ws.Range("H132").Value = 3699.5
$ws.Range("I132").Value = 1905.4286
$ws.Range("K132").Value = 5716.2858
$ws.Range("M132").Value = -3186.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 56994.332
$ws.Range("J88").Value = 56994.332
$ws.Range("L88").Value = 56994.332
$ws.Range("N88").Value = -57806.332
$ws.Range("H91").Value = 56994.332
$ws.Range("J91").Value = 56994.332
$ws.Range("L91").Value = 56994.332
$ws.Range("N91").Value = -59802.332
$ws.Range("H94").Value = 4058.6487
$ws.Range("I94").Value = 4725.7
$ws.Range("K94").Value = 4725.7
$ws.Range("M94").Value = -4274.7
$ws.Range("H105").Value = 18529.111
$ws.Range("I105").Value = 30443.5
$ws.Range("J105").Value = 8997.6
$ws.Range("K105").Value = 30443.5
$ws.Range("L105").Value = 8997.6
$ws.Range("M105").Value = -28696.5
$ws.Range("N105").Value = -12491.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 802.25
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 802.25
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H58").Value = 2148.6667
$ws.Range("J58").Value = 4014
$ws.Range("L58").Value = 4014
$ws.Range("N58").Value = -4420
$ws.Range("H99").Value = 4718.3335
$ws.Range("J99").Value = 5077.5
$ws.Range("L99").Value = 5077.5
$ws.Range("N99").Value = -8073.5
$ws.Range("H126").Value = 4718.3335
$ws.Range("J126").Value = 5077.5
$ws.Range("L126").Value = 15232.5
$ws.Range("N126").Value = -20172.5
$ws.Range("H132").Value = 4209.6113
$ws.Range("I132").Value = 4271.9287
$ws.Range("K132").Value = 12815.7861
$ws.Range("M132").Value = -10285.7861
$ws.Range("H136").Value = 2148.6667
$ws.Range("J136").Value = 4014
$ws.Range("L136").Value = 12042
$ws.Range("N136").Value = -17142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 672.2
$ws.Range("I8").Value = 672.2
$ws.Range("M8").Value = -1877.6
$ws.Range("H22").Value = 230700.1
$ws.Range("I22").Value = 11400.2
$ws.Range("J22").Value = 450000
$ws.Range("K22").Value = 34200.60000000001
$ws.Range("L22").Value = 1350000
$ws.Range("M22").Value = -34031.60000000001
$ws.Range("N22").Value = -1350338
$ws.Range("H27").Value = 230700.1
$ws.Range("I27").Value = 11400.2
$ws.Range("J27").Value = 450000
$ws.Range("K27").Value = 34200.60000000001
$ws.Range("L27").Value = 1350000
$ws.Range("M27").Value = -34098.60000000001
$ws.Range("N27").Value = -1350204
$ws.Range("H113").Value = 810.6087
$ws.Range("J113").Value = 903.8333
$ws.Range("L113").Value = 2711.4999
$ws.Range("N113").Value = -7051.4999
$ws.Range("H122").Value = 1898.6364
$ws.Range("I122").Value = 1110.6666
$ws.Range("J122").Value = 2194.125
$ws.Range("K122").Value = 9995.999400000001
$ws.Range("L122").Value = 19747.125
$ws.Range("M122").Value = -7545.999400000001
$ws.Range("N122").Value = -24647.125
$ws.Range("H132").Value = 1463.7059
$ws.Range("I132").Value = 1216.1111
$ws.Range("K132").Value = 10944.9999
$ws.Range("M132").Value = -8414.999900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 10035400
$ws.Range("J20").Value = 48333.332
$ws.Range("L20").Value = 48333.332
$ws.Range("N20").Value = -48823.332
$ws.Range("H29").Value = 9999
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H122").Value = 4101.1304
$ws.Range("I122").Value = 2113.4167
$ws.Range("J122").Value = 6269.5454
$ws.Range("K122").Value = 6340.250100000001
$ws.Range("L122").Value = 18808.6362
$ws.Range("M122").Value = -3890.250100000001
$ws.Range("N122").Value = -23708.6362
$ws.Range("H126").Value = 1691.8572
$ws.Range("I126").Value = 1620.1818
$ws.Range("J126").Value = 1954.6666
$ws.Range("K126").Value = 4860.5454
$ws.Range("L126").Value = 5863.9998
$ws.Range("M126").Value = -2390.5454
$ws.Range("N126").Value = -10803.9998
$ws.Range("H132").Value = 3297.121
$ws.Range("I132").Value = 2954.423
$ws.Range("K132").Value = 8863.269
$ws.Range("M132").Value = -6333.269

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 1016263.7
$ws.Range("I23").Value = 735196.1
$ws.Range("K23").Value = 735196.1
$ws.Range("M23").Value = -734966.1
$ws.Range("H43").Value = 5597059
$ws.Range("I43").Value = 250000
$ws.Range("J43").Value = 7242307.5
$ws.Range("K43").Value = 250000
$ws.Range("L43").Value = 7242307.5
$ws.Range("M43").Value = -249807
$ws.Range("N43").Value = -7242693.5
$ws.Range("H132").Value = 9663.362999999999
$ws.Range("I132").Value = 3521
$ws.Range("K132").Value = 10563
$ws.Range("M132").Value = -8033
$ws.Range("H136").Value = 4608.8184
$ws.Range("J136").Value = 7500
$ws.Range("L136").Value = 22500
$ws.Range("N136").Value = -27600
$ws.Range("H141").Value = 148394
$ws.Range("J141").Value = 148394
$ws.Range("L141").Value = 148394
$ws.Range("N141").Value = -158754

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 250615
$ws.Range("J105").Value = 250615
$ws.Range("L105").Value = 250615
$ws.Range("N105").Value = -257603
$ws.Range("H113").Value = 435.0625
$ws.Range("I113").Value = 260.9
$ws.Range("K113").Value = 782.6999999999999
$ws.Range("M113").Value = 1387.3
$ws.Range("H132").Value = 35061
$ws.Range("I132").Value = 45698.13
$ws.Range("K132").Value = 137094.39
$ws.Range("M132").Value = -134564.39
$ws.Range("H136").Value = 30589.324
$ws.Range("I136").Value = 32251.281
$ws.Range("K136").Value = 96753.84299999999
$ws.Range("M136").Value = -94203.84299999999
